# Box Plot Updates, Color Updates Main Figures
#
# Nudges the (x,y) position of several box-plot label textboxes that live
# inside the single group shape on slide 1. Target values below are
# expressed in points but were chosen so that the PowerPoint COM layer's
# point -> EMU conversion reproduces the exact EMU offsets from the
# authored OOXML (the runtime narrows the point value to single
# precision before multiplying by 12700 and truncating, so a handful of
# extra fractional digits are included to land on the correct EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# name -> new Left/Top (points)
$moves = @{
    "tx9"  = @{ Left = 355.8044891889764;  Top = 193.90653543307087 }
    "tx10" = @{ Left = 383.7876434952756;  Top = 218.0312598425197  }
    "tx11" = @{ Left = 418.56803899606297; Top = 249.92968503937007 }
    "tx12" = @{ Left = 465.0524445448819;  Top = 274.0544094488189  }
    "tx13" = @{ Left = 506.67063912125985; Top = 300.4351181102362  }
    "tx14" = @{ Left = 508.07937007874017; Top = 327.75244094488187 }
    "tx15" = @{ Left = 371.63472440944884; Top = 378.2068634937008  }
    "tx16" = @{ Left = 380.64149606299213; Top = 405.5241852283465  }
    "tx17" = @{ Left = 319.730157480315;   Top = 246.72094488188975 }
    "tx18" = @{ Left = 326.8405609811023;  Top = 274.0382677165354  }
}

function Apply-Moves($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        if ($moves.ContainsKey($sh.Name)) {
            $m = $moves[$sh.Name]
            $sh.Left = $m.Left
            $sh.Top = $m.Top
        }

        if ($sh.Type -eq 6) {
            Apply-Moves($sh.GroupItems)
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    Apply-Moves($slide.Shapes)
}
